$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 36
$ws1.Range("F3").Value = 183
$ws1.Range("F5").Value = 31
$ws1.Range("F6").Value = 564
$ws1.Range("F7").Value = 1761
$ws1.Range("F11").Value = 1976
$ws1.Range("F13").Value = 748
$ws1.Range("F14").Value = 452
$ws1.Range("F15").Value = 13
$ws1.Range("F16").Value = 282
$ws1.Range("F19").Value = 22
$ws1.Range("F23").Value = 4
$ws1.Range("F24").Value = 1058
$ws1.Range("F29").Value = 303

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 36
$ws4.Range("F3").Value = 183
$ws4.Range("F5").Value = 31
$ws4.Range("F6").Value = 565
$ws4.Range("F7").Value = 1761
$ws4.Range("F12").Value = 1976
$ws4.Range("F14").Value = 748
$ws4.Range("F15").Value = 452
$ws4.Range("F16").Value = 13
$ws4.Range("F17").Value = 282
$ws4.Range("F20").Value = 22
$ws4.Range("F24").Value = 4
$ws4.Range("F25").Value = 1058
$ws4.Range("F30").Value = 303
